# Update "想去人数" (want-to-go count) values in column F across sheets
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 205
$ws.Range("F4").Value = 388
$ws.Range("F5").Value = 181
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 546
$ws.Range("F8").Value = 53
$ws.Range("F9").Value = 9751
$ws.Range("F11").Value = 2632
$ws.Range("F13").Value = 2387
$ws.Range("F14").Value = 2644
$ws.Range("F16").Value = 275
$ws.Range("F17").Value = 2070
$ws.Range("F19").Value = 78
$ws.Range("F20").Value = 365
$ws.Range("F22").Value = 67
$ws.Range("F23").Value = 296
$ws.Range("F24").Value = 59
$ws.Range("F25").Value = 143
$ws.Range("F27").Value = 1282
$ws.Range("F28").Value = 1240
$ws.Range("F29").Value = 94
$ws.Range("F32").Value = 1671
$ws.Range("F33").Value = 2789
$ws.Range("F35").Value = 985
$ws.Range("F36").Value = 354
$ws.Range("F39").Value = 43
$ws.Range("F40").Value = 47
$ws.Range("F41").Value = 49
$ws.Range("F42").Value = 27
$ws.Range("F43").Value = 26

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 722
$ws.Range("F3").Value = 947
$ws.Range("F4").Value = 114
$ws.Range("F5").Value = 948

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 205
$ws.Range("F3").Value = 722
$ws.Range("F4").Value = 947
$ws.Range("F5").Value = 114
$ws.Range("F6").Value = 388
$ws.Range("F9").Value = 181
$ws.Range("F10").Value = 546
$ws.Range("F11").Value = 53
$ws.Range("F12").Value = 9751
$ws.Range("F16").Value = 2632
$ws.Range("F18").Value = 2387
$ws.Range("F19").Value = 2644
$ws.Range("F20").Value = 275
$ws.Range("F21").Value = 2070
$ws.Range("F23").Value = 78
$ws.Range("F24").Value = 365
$ws.Range("F26").Value = 296
$ws.Range("F27").Value = 59
$ws.Range("F28").Value = 143
$ws.Range("F30").Value = 1282
$ws.Range("F31").Value = 1240
$ws.Range("F32").Value = 94
$ws.Range("F35").Value = 1671
$ws.Range("F37").Value = 2789
$ws.Range("F38").Value = 985
$ws.Range("F41").Value = 354
$ws.Range("F46").Value = 49
$ws.Range("F47").Value = 27
$ws.Range("F48").Value = 26

